$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Price cells whose new values look like plain numbers need to be forced to
# Text format first, since the source data stores prices as text strings
# (e.g. "6.75"), not numeric values.
$textCells = @("D5","D6","D9","D12","D13","D17","D19","D20","D22","D23","D24","D25","D28","D29","D31","D32","D34","D35","D36","D37","D38","D41","D42","D43","D44","D46","D47","D48","D51")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = "63.552.71"
$ws.Range("E2").Value = "  -1.00%  "
$ws.Range("D3").Value = "3.109.97"
$ws.Range("E3").Value = "  +0.99%  "
$ws.Range("E4").Value = "  -0.07%  "
$ws.Range("D5").Value = "558.48"
$ws.Range("E5").Value = "  +1.02%  "
$ws.Range("D6").Value = "139.06"
$ws.Range("E6").Value = "  -2.68%  "
$ws.Range("E7").Value = "  +0.06%  "
$ws.Range("D8").Value = "3.100.64"
$ws.Range("E8").Value = "  +1.01%  "
$ws.Range("D9").Value = "0.499"
$ws.Range("E9").Value = "  +1.56%  "
$ws.Range("E10").Value = "  +3.47%  "
$ws.Range("E11").Value = "  +4.85%  "
$ws.Range("D12").Value = "0.459"
$ws.Range("E12").Value = "  +1.41%  "
$ws.Range("D13").Value = "35.59"
$ws.Range("E13").Value = "  -0.83%  "
$ws.Range("E14").Value = "  +0.53%  "
$ws.Range("D15").Value = "3.612.61"
$ws.Range("E15").Value = "  +1.00%  "
$ws.Range("D16").Value = "63.587.60"
$ws.Range("E16").Value = "  -1.15%  "
$ws.Range("D17").Value = "0.112"
$ws.Range("E17").Value = "  +0.39%  "
$ws.Range("D18").Value = "3.110.42"
$ws.Range("E18").Value = "  +0.98%  "
$ws.Range("D19").Value = "507.06"
$ws.Range("E19").Value = "  +3.81%  "
$ws.Range("D20").Value = "6.75"
$ws.Range("E21").Value = "  +0.46%  "
$ws.Range("D22").Value = "0.714"
$ws.Range("E22").Value = "  +3.94%  "
$ws.Range("D23").Value = "7.39"
$ws.Range("E23").Value = "  +2.51%  "
$ws.Range("B24").Value = "Litecoin"
$ws.Range("C24").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D24").Value = "78.59"
$ws.Range("E24").Value = "  +0.86%  "
$ws.Range("B25").Value = "InternetComputer(DFINITY)"
$ws.Range("C25").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D25").Value = "12.49"
$ws.Range("E25").Value = "  +0.29%  "
$ws.Range("E26").Value = "  -0.06%  "
$ws.Range("E27").Value = "  +2.63%  "
$ws.Range("D28").Value = "8.32"
$ws.Range("E28").Value = "  -0.07%  "
$ws.Range("D29").Value = "2.06"
$ws.Range("E29").Value = "  -1.19%  "
$ws.Range("E30").Value = "  -0.04%  "
$ws.Range("D31").Value = "26.37"
$ws.Range("E31").Value = "  +2.28%  "
$ws.Range("D32").Value = "2.55"
$ws.Range("E32").Value = "  -3.44%  "
$ws.Range("E33").Value = "  -0.93%  "
$ws.Range("B34").Value = "OKB"
$ws.Range("C34").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D34").Value = "59.23"
$ws.Range("E34").Value = "  +13.70%  "
$ws.Range("B35").Value = "Bittensor"
$ws.Range("C35").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D35").Value = "538.30"
$ws.Range("E35").Value = "  -7.97%  "
$ws.Range("D36").Value = "5.95"
$ws.Range("E36").Value = "  +0.34%  "
$ws.Range("D37").Value = "5.24"
$ws.Range("E37").Value = "  -2.51%  "
$ws.Range("D38").Value = "0.0419"
$ws.Range("E38").Value = "  +4.05%  "
$ws.Range("E39").Value = "  +1.64%  "
$ws.Range("D40").Value = "3.095.85"
$ws.Range("E40").Value = "  +3.26%  "
$ws.Range("D41").Value = "0.119"
$ws.Range("E41").Value = "  +1.26%  "
$ws.Range("D42").Value = "8.18"
$ws.Range("D43").Value = "2.69"
$ws.Range("E43").Value = "  -5.52%  "
$ws.Range("D44").Value = "0.259"
$ws.Range("E44").Value = "  +5.69%  "
$ws.Range("E45").Value = "  +0.05%  "
$ws.Range("D46").Value = "2.12"
$ws.Range("E46").Value = "  +0.91%  "
$ws.Range("D47").Value = "121.18"
$ws.Range("E47").Value = "  +1.56%  "
$ws.Range("D48").Value = "24.31"
$ws.Range("E48").Value = "  -3.37%  "
$ws.Range("E49").Value = "  -0.05%  "
$ws.Range("D50").Value = "0.0₃0503"
$ws.Range("E50").Value = "  -5.19%  "
$ws.Range("D51").Value = "2.36"
$ws.Range("E51").Value = "  +62.31%  "

# Restore default styling on the cells we temporarily reformatted so the
# only observable change is the cell content, not its style index.
foreach ($addr in $textCells) {
    $ws.Range($addr).Style = "Normal"
}
